$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# New sheet: 仪表盘 (gauge chart) - appended after the last sheet
# ---------------------------------------------------------------
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGauge = $wb.Worksheets.Add($null, $after)
$wsGauge.Name = "仪表盘"
$wsGauge.Range("A1").Value = "完成率"
$wsGauge.Range("B1").Value = 55.5
$wsGauge.Range("B3").Select()

# ---------------------------------------------------------------
# New sheet: 漏斗图 (funnel chart) - appended after 仪表盘
# ---------------------------------------------------------------
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFunnel = $wb.Worksheets.Add($null, $after)
$wsFunnel.Name = "漏斗图"
$wsFunnel.Range("A1").Value = "展现"
$wsFunnel.Range("B1").Value = 100
$wsFunnel.Range("A2").Value = "点击"
$wsFunnel.Range("B2").Value = 80
$wsFunnel.Range("A3").Value = "访问"
$wsFunnel.Range("B3").Value = 60
$wsFunnel.Range("A4").Value = "咨询"
$wsFunnel.Range("B4").Value = 40
$wsFunnel.Range("A5").Value = "订单"
$wsFunnel.Range("B5").Value = 20
$wsFunnel.Range("F43").Select()

# ---------------------------------------------------------------
# New sheet: 雷达图 (radar chart) - appended after 漏斗图
# Data rows written first, header labels (row 1) written last,
# matching the original authoring order.
# ---------------------------------------------------------------
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRadar = $wb.Worksheets.Add($null, $after)
$wsRadar.Name = "雷达图"

$wsRadar.Range("A2").Value = "销售"
$wsRadar.Range("B2").Value = 6500
$wsRadar.Range("C2").Value = 4300
$wsRadar.Range("D2").Value = 5000

$wsRadar.Range("A3").Value = "经理"
$wsRadar.Range("B3").Value = 16000
$wsRadar.Range("C3").Value = 10000
$wsRadar.Range("D3").Value = 14000

$wsRadar.Range("A4").Value = "信息技术"
$wsRadar.Range("B4").Value = 30000
$wsRadar.Range("C4").Value = 28000
$wsRadar.Range("D4").Value = 28000

$wsRadar.Range("A5").Value = "客服"
$wsRadar.Range("B5").Value = 38000
$wsRadar.Range("C5").Value = 35000
$wsRadar.Range("D5").Value = 31000

$wsRadar.Range("A6").Value = "研发"
$wsRadar.Range("B6").Value = 52000
$wsRadar.Range("C6").Value = 50000
$wsRadar.Range("D6").Value = 42000

$wsRadar.Range("A7").Value = "市场"
$wsRadar.Range("B7").Value = 25000
$wsRadar.Range("C7").Value = 19000
$wsRadar.Range("D7").Value = 21000

$wsRadar.Range("C1").Value = "预算分配"
$wsRadar.Range("D1").Value = "实际开销"

$wsRadar.Range("F18").Select()
